$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 54; everything from old row 54 downward shifts down by one
# (old row 54 -> new row 55, ..., old row 74 -> new row 75).
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly price-report record.
$ws.Cells.Item(54,1).Value = 9
$ws.Cells.Item(54,2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(54,3).Value = 'Metropolitana'
$ws.Cells.Item(54,4).Value = '2021-11-10'
$ws.Cells.Item(54,5).Value = 13
$ws.Cells.Item(54,6).Value = 100112022
$ws.Cells.Item(54,7).Value = 'Arveja Verde'
$ws.Cells.Item(54,8).Value = 'Sin especificar'
$ws.Cells.Item(54,9).Value = 'Primera'
$ws.Cells.Item(54,10).Value = 34
$ws.Cells.Item(54,11).Value = 14000
$ws.Cells.Item(54,12).Value = 15000
$ws.Cells.Item(54,13).Value = 14500
$ws.Cells.Item(54,14).Value = '$/saco 25 kilos'
$ws.Cells.Item(54,15).Value = 'Región Metropolitana'
$ws.Cells.Item(54,16).Value = 580
$ws.Cells.Item(54,17).Value = 25
$ws.Cells.Item(54,18).Value = 'Hortaliza'
